$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data (row 2), pushing the existing
# records (old rows 2-20) down to rows 5-23.
$ws.Rows("2:4").Insert()

# The freshly inserted rows inherit the header row's bold/centered style;
# clear that so they look like ordinary data rows.
$ws.Range("A2:T4").ClearFormats()

# Column D carries a date-formatted style (same as the rest of the table) -
# copy it from the row right below (now row 5) onto the 3 new date cells.
$ws.Range("D5").Copy()
$ws.Range("D2:D4").PasteSpecial(-4122)

# New daily price records for 2022-03-11 (serial 44631), one row per
# "Calidad" grade, matching the layout of the existing rows.
$ws.Range("A2").Value2 = 8
$ws.Range("B2").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value2 = "Coquimbo"
$ws.Range("D2").Value2 = 44631
$ws.Range("E2").Value2 = 4
$ws.Range("F2").Value2 = "Fruta"
$ws.Range("G2").Value2 = 100107
$ws.Range("H2").Value2 = "Otros"
$ws.Range("I2").Value2 = 100107011
$ws.Range("J2").Value2 = "Tuna"
$ws.Range("K2").Value2 = "Sin especificar"
$ws.Range("L2").Value2 = "Especial"
$ws.Range("M2").Value2 = 240
$ws.Range("N2").Value2 = 15000
$ws.Range("O2").Value2 = 16000
$ws.Range("P2").Value2 = 15500
$ws.Range("Q2").Value2 = "$/caja 18 kilos"
$ws.Range("R2").Value2 = "Provincia de Limarí"
$ws.Range("S2").Value2 = 861
$ws.Range("T2").Value2 = 18

$ws.Range("A3").Value2 = 8
$ws.Range("B3").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value2 = "Coquimbo"
$ws.Range("D3").Value2 = 44631
$ws.Range("E3").Value2 = 4
$ws.Range("F3").Value2 = "Fruta"
$ws.Range("G3").Value2 = 100107
$ws.Range("H3").Value2 = "Otros"
$ws.Range("I3").Value2 = 100107011
$ws.Range("J3").Value2 = "Tuna"
$ws.Range("K3").Value2 = "Sin especificar"
$ws.Range("L3").Value2 = "Primera"
$ws.Range("M3").Value2 = 248
$ws.Range("N3").Value2 = 12000
$ws.Range("O3").Value2 = 13000
$ws.Range("P3").Value2 = 12516
$ws.Range("Q3").Value2 = "$/caja 18 kilos"
$ws.Range("R3").Value2 = "Provincia de Limarí"
$ws.Range("S3").Value2 = 695
$ws.Range("T3").Value2 = 18

$ws.Range("A4").Value2 = 8
$ws.Range("B4").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value2 = "Coquimbo"
$ws.Range("D4").Value2 = 44631
$ws.Range("E4").Value2 = 4
$ws.Range("F4").Value2 = "Fruta"
$ws.Range("G4").Value2 = 100107
$ws.Range("H4").Value2 = "Otros"
$ws.Range("I4").Value2 = 100107011
$ws.Range("J4").Value2 = "Tuna"
$ws.Range("K4").Value2 = "Sin especificar"
$ws.Range("L4").Value2 = "Segunda"
$ws.Range("M4").Value2 = 200
$ws.Range("N4").Value2 = 9000
$ws.Range("O4").Value2 = 10000
$ws.Range("P4").Value2 = 9500
$ws.Range("Q4").Value2 = "$/caja 18 kilos"
$ws.Range("R4").Value2 = "Provincia de Limarí"
$ws.Range("S4").Value2 = 528
$ws.Range("T4").Value2 = 18
